$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(32)
$rng = $p.Range
$rng.InsertParagraphAfter()
$newp = $d.Paragraphs.Item(33)
Write-Output ("newp text=[" + $newp.Range.Text + "]")
$newp.Range.Text = "Section number changed"
Write-Output ("newp text now=[" + $newp.Range.Text + "]")
# delete old paragraph 32 (now still index 32)
$p2 = $d.Paragraphs.Item(32)
Write-Output ("p2 text=[" + $p2.Range.Text + "]")
$p2.Range.Delete()
